$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values are plain numeric-looking strings (e.g. "578.56").
# The source data stores these as text (prices formatted with literal dot thousands
# separators elsewhere in the column), so force a Text number format before assigning
# the value to prevent Excel from auto-converting them to numbers.
$textCells = @("D5", "D6", "D10", "D12", "D13", "D16", "D19", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D40", "D41", "D45", "D46", "D47", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (column D) and Volume(1h) (column E) values scraped by the
# GitHub Actions job, row by row.
$ws.Range("D2").Value = "69.872.77"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.568.32"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "578.56"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").Value = "188.46"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D8").Value = "3.564.95"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "0.176"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "55.82"
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("D13").Value = "0.0000301"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").Value = "4.144.89"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "19.78"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "3.573.74"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "69.699.82"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "12.60"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "474.72"
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").Value = "19.22"
$ws.Range("E23").Value = "  +14.52%  "
$ws.Range("D24").Value = "5.07"
$ws.Range("E24").Value = "  -7.89%  "
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "93.49"
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "10.98"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "9.30"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "32.21"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "7.73"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").Value = "12.17"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "66.20"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "579.57"
$ws.Range("E35").Value = "  -5.88%  "
$ws.Range("D36").Value = "39.00"
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "3.25"
$ws.Range("E40").Value = "  +18.87%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.139"
$ws.Range("E41").Value = "  -5.95%  "
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("D43").Value = "3.226.37"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("E44").Value = "  +6.87%  "
$ws.Range("D45").Value = "3.07"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").Value = "0.0441"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.37"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "9.48"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -5.90%  "
